# Add a new "isaterms" column (E) to Sheet1 of the mapping_file workbook.
# Each row gets the corresponding isa-terms-style (snake_case) mapping value
# for the ISA class named in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for the new column
$ws.Range("E1").Value = "isaterms"

# Mapping of row number -> isaterms value, in the same row order as the
# existing table (rows 2-23, one per ISA class already present in columns A-D).
# Row 8 (material_entity) is filled in last, after all the others, which
# matches the order new shared strings were appended in the original edit.
$orderedRows = 2, 3, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 8

$values = @{
    2  = "assay_type"
    3  = "comment"
    4  = "file"
    5  = "factor"
    6  = "factor_value"
    7  = "investigation"
    8  = "material_entity"
    9  = "characteristic"
    10 = "characteristic_value"
    11 = "NA"
    12 = "ontology_annotation"
    13 = "ontology"
    14 = "organization"
    15 = "person"
    16 = "process"
    17 = "parameter_value"
    18 = "protocol"
    19 = "protocol_parameter"
    20 = "publication"
    21 = "sample"
    22 = "source"
    23 = "study"
}

foreach ($row in $orderedRows) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}

# Update the view to match: scroll so row 2 is at the top, and select C17
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
